$wb = $excel.ActiveWorkbook

# Column F ("想去人数" / want-to-go count) updates that occurred when the
# site data was regenerated. Same updates apply identically to the
# "展览" (Exhibition) sheet and the "全部类型" (All types) sheet, which
# mirror each other's data.

$updates = @{
    "F3"  = 372
    "F5"  = 411
    "F7"  = 121
    "F9"  = 58
    "F11" = 117
    "F12" = 1136
    "F13" = 1465
    "F14" = 324
    "F15" = 361
    "F17" = 99
    "F18" = 159
    "F21" = 257
    "F22" = 278
    "F23" = 303
    "F24" = 1677
    "F25" = 61
    "F27" = 165
    "F28" = 632
    "F30" = 75
    "F31" = 3967
    "F33" = 464
    "F34" = 238
    "F35" = 1014
    "F36" = 104
    "F37" = 46
    "F39" = 104
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($cellRef in $updates.Keys) {
        $ws.Range($cellRef).Value = $updates[$cellRef]
    }
}

$wb.Save()
